# Update cryptocurrency price/volume data in the worksheet.
# Values that could be misinterpreted as numbers by Excel are written with a
# leading apostrophe so they remain plain text (matching the workbook's
# original inlineStr cell representation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '53.877.54'
$ws.Range("E2").Value = '  -4.21%  '
$ws.Range("D3").Value = '2.223.79'
$ws.Range("E3").Value = '  -6.11%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''490.04'
$ws.Range("E5").Value = '  -3.18%  '
$ws.Range("D6").Value = '''126.62'
$ws.Range("E6").Value = '  -2.35%  '
$ws.Range("E7").Value = '  -0.42%  '
$ws.Range("E8").Value = '  -3.40%  '
$ws.Range("D9").Value = '2.259.26'
$ws.Range("E9").Value = '  -4.92%  '
$ws.Range("D10").Value = '''0.0929'
$ws.Range("E10").Value = '  -5.54%  '
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '''0.321'
$ws.Range("E12").Value = '  -2.01%  '
$ws.Range("D13").Value = '''4.63'
$ws.Range("E13").Value = '  -4.47%  '
$ws.Range("D14").Value = '2.623.25'
$ws.Range("E14").Value = '  -6.05%  '
$ws.Range("D15").Value = '''21.44'
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("D16").Value = '53.822.83'
$ws.Range("E16").Value = '  -4.30%  '
$ws.Range("E17").Value = '  -2.89%  '
$ws.Range("D18").Value = '2.234.93'
$ws.Range("E18").Value = '  -5.87%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '''4.03'
$ws.Range("E19").Value = '  -0.14%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '''9.74'
$ws.Range("E20").Value = '  -2.72%  '
$ws.Range("D21").Value = '''297.64'
$ws.Range("E21").Value = '  -3.57%  '
$ws.Range("D22").Value = '''6.26'
$ws.Range("E22").Value = '  -0.38%  '
$ws.Range("D23").Value = '''0.995'
$ws.Range("E23").Value = '  -0.38%  '
$ws.Range("D24").Value = '''63.70'
$ws.Range("E24").Value = '  -3.54%  '
$ws.Range("D25").Value = '''0.996'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").Value = '''0.372'
$ws.Range("E26").Value = '  +0.37%  '
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '''0.147'
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '2.336.97'
$ws.Range("E28").Value = '  -5.94%  '
$ws.Range("D29").Value = '''7.09'
$ws.Range("E29").Value = '  -1.57%  '
$ws.Range("D30").Value = '''162.70'
$ws.Range("E30").Value = '  -6.00%  '
$ws.Range("E31").Value = '  -3.27%  '
$ws.Range("D32").Value = '0.0₃0675'
$ws.Range("E32").Value = '  -4.85%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").Value = '''0.997'
$ws.Range("E33").Value = '  -0.24%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").Value = '''5.80'
$ws.Range("E34").Value = '  -0.37%  '
$ws.Range("E35").Value = '  -0.52%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("D37").Value = '''17.41'
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("E38").Value = '  +0.85%  '
$ws.Range("D39").Value = '''0.836'
$ws.Range("E39").Value = '  +1.84%  '
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("D41").Value = '''35.34'
$ws.Range("E41").Value = '  -2.66%  '
$ws.Range("D42").Value = '''0.374'
$ws.Range("E42").Value = '  +0.51%  '
$ws.Range("E43").Value = '  +0.49%  '
$ws.Range("E44").Value = '  -1.52%  '
$ws.Range("D45").Value = '''126.87'
$ws.Range("E45").Value = '  +0.84%  '
$ws.Range("D46").Value = '''4.87'
$ws.Range("E46").Value = '  +1.74%  '
$ws.Range("D48").Value = '''242.47'
$ws.Range("E48").Value = '  +1.81%  '
$ws.Range("D49").Value = '''0.542'
$ws.Range("E49").Value = '  -3.57%  '
$ws.Range("D50").Value = '''0.0476'
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("E51").Value = '  -2.00%  '
